$d = $word.ActiveDocument

# Update the date title
$d.Paragraphs.Item(1).Range.Text = "2025-07-22 Tuesday"

# Update the answer table (5 columns x 20 rows)
$tbl = $d.Tables.Item(1)

$tbl.Cell(1, 1).Range.Text = "9+84=93"
$tbl.Cell(1, 2).Range.Text = "79+1=80"
$tbl.Cell(1, 3).Range.Text = "56+37=93"
$tbl.Cell(1, 4).Range.Text = "67-38=29"
$tbl.Cell(1, 5).Range.Text = "31+41=72"
$tbl.Cell(2, 1).Range.Text = "78+9=87"
$tbl.Cell(2, 2).Range.Text = "68-29=39"
$tbl.Cell(2, 3).Range.Text = "86+1=87"
$tbl.Cell(2, 4).Range.Text = "10+76=86"
$tbl.Cell(2, 5).Range.Text = "58-34=24"
$tbl.Cell(3, 1).Range.Text = "97-80=17"
$tbl.Cell(3, 2).Range.Text = "83-5=78"
$tbl.Cell(3, 3).Range.Text = "68-46=22"
$tbl.Cell(3, 4).Range.Text = "72-71=1"
$tbl.Cell(3, 5).Range.Text = "99-55=44"
$tbl.Cell(4, 1).Range.Text = "63-42=21"
$tbl.Cell(4, 2).Range.Text = "20+67=87"
$tbl.Cell(4, 3).Range.Text = "48+40=88"
$tbl.Cell(4, 4).Range.Text = "72+25=97"
$tbl.Cell(4, 5).Range.Text = "79+11=90"
$tbl.Cell(5, 1).Range.Text = "76-16=60"
$tbl.Cell(5, 2).Range.Text = "16-5=11"
$tbl.Cell(5, 3).Range.Text = "4+30=34"
$tbl.Cell(5, 4).Range.Text = "12-8=4"
$tbl.Cell(5, 5).Range.Text = "28+51=79"
$tbl.Cell(6, 1).Range.Text = "42-27=15"
$tbl.Cell(6, 2).Range.Text = "87-2=85"
$tbl.Cell(6, 3).Range.Text = "74-57=17"
$tbl.Cell(6, 4).Range.Text = "19+51=70"
$tbl.Cell(6, 5).Range.Text = "21+6=27"
$tbl.Cell(7, 1).Range.Text = "39+29=68"
$tbl.Cell(7, 2).Range.Text = "99-15=84"
$tbl.Cell(7, 3).Range.Text = "2+92=94"
$tbl.Cell(7, 4).Range.Text = "3+62=65"
$tbl.Cell(7, 5).Range.Text = "82+17=99"
$tbl.Cell(8, 1).Range.Text = "27+69=96"
$tbl.Cell(8, 2).Range.Text = "0+94=94"
$tbl.Cell(8, 3).Range.Text = "66+31=97"
$tbl.Cell(8, 4).Range.Text = "87+11=98"
$tbl.Cell(8, 5).Range.Text = "42+33=75"
$tbl.Cell(9, 1).Range.Text = "89-13=76"
$tbl.Cell(9, 2).Range.Text = "2+95=97"
$tbl.Cell(9, 3).Range.Text = "37+12=49"
$tbl.Cell(9, 4).Range.Text = "70-3=67"
$tbl.Cell(9, 5).Range.Text = "29+51=80"
$tbl.Cell(10, 1).Range.Text = "52-15=37"
$tbl.Cell(10, 2).Range.Text = "67-46=21"
$tbl.Cell(10, 3).Range.Text = "6+38=44"
$tbl.Cell(10, 4).Range.Text = "2+2=4"
$tbl.Cell(10, 5).Range.Text = "31+40=71"
$tbl.Cell(11, 1).Range.Text = "40+27=67"
$tbl.Cell(11, 2).Range.Text = "48+2=50"
$tbl.Cell(11, 3).Range.Text = "15+76=91"
$tbl.Cell(11, 4).Range.Text = "82-5=77"
$tbl.Cell(11, 5).Range.Text = "56-33=23"
$tbl.Cell(12, 1).Range.Text = "13+7=20"
$tbl.Cell(12, 2).Range.Text = "72-4=68"
$tbl.Cell(12, 3).Range.Text = "40+52=92"
$tbl.Cell(12, 4).Range.Text = "79-0=79"
$tbl.Cell(12, 5).Range.Text = "33-6=27"
$tbl.Cell(13, 1).Range.Text = "72-43=29"
$tbl.Cell(13, 2).Range.Text = "18+49=67"
$tbl.Cell(13, 3).Range.Text = "50+40=90"
$tbl.Cell(13, 4).Range.Text = "21-6=15"
$tbl.Cell(13, 5).Range.Text = "96-57=39"
$tbl.Cell(14, 1).Range.Text = "60-5=55"
$tbl.Cell(14, 2).Range.Text = "8+20=28"
$tbl.Cell(14, 3).Range.Text = "77-53=24"
$tbl.Cell(14, 4).Range.Text = "65+22=87"
$tbl.Cell(14, 5).Range.Text = "60-23=37"
$tbl.Cell(15, 1).Range.Text = "83-79=4"
$tbl.Cell(15, 2).Range.Text = "70-1=69"
$tbl.Cell(15, 3).Range.Text = "13+51=64"
$tbl.Cell(15, 4).Range.Text = "49-21=28"
$tbl.Cell(15, 5).Range.Text = "40+56=96"
$tbl.Cell(16, 1).Range.Text = "87-86=1"
$tbl.Cell(16, 2).Range.Text = "2+1=3"
$tbl.Cell(16, 3).Range.Text = "41-40=1"
$tbl.Cell(16, 4).Range.Text = "89-30=59"
$tbl.Cell(16, 5).Range.Text = "72-62=10"
$tbl.Cell(17, 1).Range.Text = "98-51=47"
$tbl.Cell(17, 2).Range.Text = "41-23=18"
$tbl.Cell(17, 3).Range.Text = "54+38=92"
$tbl.Cell(17, 4).Range.Text = "6+43=49"
$tbl.Cell(17, 5).Range.Text = "48-29=19"
$tbl.Cell(18, 1).Range.Text = "80-73=7"
$tbl.Cell(18, 2).Range.Text = "4+72=76"
$tbl.Cell(18, 3).Range.Text = "93-4=89"
$tbl.Cell(18, 4).Range.Text = "83-17=66"
$tbl.Cell(18, 5).Range.Text = "71+26=97"
$tbl.Cell(19, 1).Range.Text = "97-97=0"
$tbl.Cell(19, 2).Range.Text = "55+3=58"
$tbl.Cell(19, 3).Range.Text = "58+1=59"
$tbl.Cell(19, 4).Range.Text = "99-43=56"
$tbl.Cell(19, 5).Range.Text = "91+4=95"
$tbl.Cell(20, 1).Range.Text = "30+23=53"
$tbl.Cell(20, 2).Range.Text = "46+16=62"
$tbl.Cell(20, 3).Range.Text = "92-90=2"
$tbl.Cell(20, 4).Range.Text = "41+14=55"
$tbl.Cell(20, 5).Range.Text = "94-14=80"
